# The canonical diff only touches PowerPoint's internal co-authoring
# revision log (ppt/revisionInfo.xml + ppt/changesInfos/changesInfo1.xml):
# it records that, on slide 267 (cId="3608403467"), a graphicFrame
# (id="2", creationId="{04FA92D2-B5C6-C5AD-C98C-7F030AF53DAE}") was
# added and then deleted again within the same editing session
# (chg="add del mod"), leaving the slide's visible content unchanged
# (the object never survives to the final saved state - compare with
# the sibling picChg id="7" on the same slide, which is also "add del
# mod" and likewise absent from the final shape tree).
#
# Reproduce that user action against the real object model: insert a
# table (a graphicFrame) onto slide 267 and immediately delete it
# again, so the slide ends up exactly as it started.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(2)

$frame = $s.Shapes.AddTable(2, 2, 2621280, 1524000, 3810000, 2540000)
$frame.Delete()
